$wb = $excel.ActiveWorkbook

# ---- Sheet ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 332.07144
$ws.Range("I2").Value = 345.75
$ws.Range("K2").Value = 345.75
$ws.Range("M2").Value = -232.75
$ws.Range("H32").Value = 19444.445
$ws.Range("J32").Value = 21000
$ws.Range("L32").Value = 21000
$ws.Range("N32").Value = -21652
$ws.Range("H38").Value = 1792.75
$ws.Range("I38").Value = 1723.8334
$ws.Range("J38").Value = 1999.5
$ws.Range("K38").Value = 5171.5002
$ws.Range("L38").Value = 5998.5
$ws.Range("M38").Value = -4799.5002
$ws.Range("N38").Value = -6742.5
$ws.Range("H87").Value = 0
$ws.Range("J87").Value = 0
$ws.Range("L87").Value = 0
$ws.Range("N87").Value = $null
$ws.Range("H90").Value = 0
$ws.Range("J90").Value = 0
$ws.Range("L90").Value = 0
$ws.Range("N90").Value = $null
$ws.Range("H108").Value = 0
$ws.Range("J108").Value = 0
$ws.Range("L108").Value = 0
$ws.Range("N108").Value = $null
$ws.Range("H138").Value = 2216.6128
$ws.Range("J138").Value = 2350.0952
$ws.Range("L138").Value = 7050.285600000001
$ws.Range("N138").Value = -17330.2856

# ---- Sheet ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H13").Value = 8000401
$ws.Range("I13").Value = 8000401
$ws.Range("K13").Value = 8000401
$ws.Range("M13").Value = -8000257
$ws.Range("H23").Value = 5000
$ws.Range("I23").Value = 5000
$ws.Range("K23").Value = 5000
$ws.Range("M23").Value = -4741
$ws.Range("H29").Value = 3632
$ws.Range("I29").Value = 3632
$ws.Range("K29").Value = 3632
$ws.Range("M29").Value = -3324
$ws.Range("H44").Value = 54000
$ws.Range("J44").Value = 54000
$ws.Range("L44").Value = 54000
$ws.Range("N44").Value = -54976
$ws.Range("H80").Value = 2000
$ws.Range("I80").Value = 2000
$ws.Range("J80").Value = 0
$ws.Range("K80").Value = 2000
$ws.Range("L80").Value = 0
$ws.Range("M80").Value = -1002
$ws.Range("N80").Value = $null
$ws.Range("H83").Value = 2000
$ws.Range("I83").Value = 2000
$ws.Range("J83").Value = 0
$ws.Range("K83").Value = 6000
$ws.Range("L83").Value = 0
$ws.Range("M83").Value = -1008
$ws.Range("N83").Value = $null
$ws.Range("H88").Value = 1605
$ws.Range("I88").Value = 1099
$ws.Range("J88").Value = 1835
$ws.Range("K88").Value = 1099
$ws.Range("L88").Value = 1835
$ws.Range("M88").Value = -693
$ws.Range("N88").Value = -2647
$ws.Range("H91").Value = 1605
$ws.Range("I91").Value = 1099
$ws.Range("J91").Value = 1835
$ws.Range("K91").Value = 1099
$ws.Range("L91").Value = 1835
$ws.Range("M91").Value = 305
$ws.Range("N91").Value = -4643

# ---- Sheet BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H35").Value = 46037
$ws.Range("J35").Value = 60074
$ws.Range("L35").Value = 60074
$ws.Range("N35").Value = -60694
$ws.Range("H86").Value = 16799.6
$ws.Range("I86").Value = 21666.166
$ws.Range("J86").Value = 9499.75
$ws.Range("K86").Value = 21666.166
$ws.Range("L86").Value = 9499.75
$ws.Range("M86").Value = -20543.166
$ws.Range("N86").Value = -11745.75
$ws.Range("H89").Value = 16799.6
$ws.Range("I89").Value = 21666.166
$ws.Range("J89").Value = 9499.75
$ws.Range("K89").Value = 108330.83
$ws.Range("L89").Value = 47498.75
$ws.Range("M89").Value = -102714.83
$ws.Range("N89").Value = -58730.75
$ws.Range("H98").Value = 120000
$ws.Range("J98").Value = 120000
$ws.Range("L98").Value = 120000
$ws.Range("N98").Value = -125990

# ---- Sheet CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 97
$ws.Range("I4").Value = 97
$ws.Range("K4").Value = 97
$ws.Range("M4").Value = 15
$ws.Range("H36").Value = 3000
$ws.Range("I36").Value = 3000
$ws.Range("K36").Value = 3000
$ws.Range("M36").Value = -2612
$ws.Range("H40").Value = 3000
$ws.Range("I40").Value = 3000
$ws.Range("K40").Value = 3000
$ws.Range("M40").Value = -2840
$ws.Range("H41").Value = 19749.5
$ws.Range("J41").Value = 0
$ws.Range("L41").Value = 0
$ws.Range("N41").Value = $null
$ws.Range("H50").Value = 30000
$ws.Range("I50").Value = 30000
$ws.Range("J50").Value = 0
$ws.Range("K50").Value = 30000
$ws.Range("L50").Value = 0
$ws.Range("M50").Value = -29375
$ws.Range("N50").Value = $null
$ws.Range("H51").Value = 38570.57
$ws.Range("J51").Value = 37995
$ws.Range("L51").Value = 37995
$ws.Range("N51").Value = -39467
$ws.Range("H60").Value = 15627.143
$ws.Range("J60").Value = 47000
$ws.Range("L60").Value = 47000
$ws.Range("N60").Value = -48022
$ws.Range("H61").Value = 38570.57
$ws.Range("J61").Value = 37995
$ws.Range("L61").Value = 37995
$ws.Range("M61").Value = 38318.5
$ws.Range("N61").Value = -38691
$ws.Range("H68").Value = 73076.92
$ws.Range("J68").Value = 75000
$ws.Range("L68").Value = 75000
$ws.Range("N68").Value = -76498
$ws.Range("H71").Value = 73076.92
$ws.Range("J71").Value = 75000
$ws.Range("L71").Value = 225000
$ws.Range("N71").Value = -232488
$ws.Range("H141").Value = 462507.2
$ws.Range("J141").Value = 462507.2
$ws.Range("L141").Value = 462507.2
$ws.Range("N141").Value = -472867.2

# ---- Sheet CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 697.375
$ws.Range("J2").Value = 1096.2
$ws.Range("L2").Value = 6577.200000000001
$ws.Range("N2").Value = -6803.200000000001
$ws.Range("H4").Value = 1317117.6
$ws.Range("I4").Value = 17950.25
$ws.Range("J4").Value = 4286643
$ws.Range("K4").Value = 53850.75
$ws.Range("L4").Value = 12859929
$ws.Range("M4").Value = -53738.75
$ws.Range("N4").Value = -12860153
$ws.Range("H7").Value = 91572.766
$ws.Range("I7").Value = 117039.6
$ws.Range("J7").Value = 6683.3335
$ws.Range("K7").Value = 351118.8
$ws.Range("L7").Value = 20050.0005
$ws.Range("M7").Value = -351006.8
$ws.Range("N7").Value = -20274.0005
$ws.Range("H11").Value = 0
$ws.Range("I11").Value = 0
$ws.Range("K11").Value = 0
$ws.Range("M11").Value = $null
$ws.Range("H12").Value = 73.90909000000001
$ws.Range("I12").Value = 71
$ws.Range("J12").Value = 74.55556
$ws.Range("K12").Value = 213
$ws.Range("L12").Value = 223.66668
$ws.Range("M12").Value = -40
$ws.Range("N12").Value = -569.66668

# ---- Sheet GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H27").Value = 0
$ws.Range("I27").Value = 0
$ws.Range("J27").Value = 0
$ws.Range("K27").Value = 0
$ws.Range("L27").Value = 0
$ws.Range("M27").Value = $null
$ws.Range("N27").Value = $null
$ws.Range("H43").Value = 8188.091
$ws.Range("J43").Value = 42999.5
$ws.Range("L43").Value = 42999.5
$ws.Range("N43").Value = -43301.5
$ws.Range("H46").Value = 37661.055
$ws.Range("I46").Value = 15999.5
$ws.Range("J46").Value = 40368.75
$ws.Range("K46").Value = 15999.5
$ws.Range("L46").Value = 40368.75
$ws.Range("M46").Value = -15843.5
$ws.Range("N46").Value = -40680.75
$ws.Range("H57").Value = 22999.6
$ws.Range("J57").Value = 35832.668
$ws.Range("L57").Value = 35832.668
$ws.Range("N57").Value = -37472.668
$ws.Range("H136").Value = 40956.168
$ws.Range("J136").Value = 31600.646
$ws.Range("L136").Value = 94801.93799999999
$ws.Range("N136").Value = -99901.93799999999

# ---- Sheet LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1490
$ws.Range("J22").Value = 1492.1666
$ws.Range("L22").Value = 1492.1666
$ws.Range("N22").Value = -2082.1666
$ws.Range("H27").Value = 1490
$ws.Range("J27").Value = 1492.1666
$ws.Range("L27").Value = 1492.1666
$ws.Range("N27").Value = -1706.1666
$ws.Range("H46").Value = 4966
$ws.Range("I46").Value = 4900
$ws.Range("K46").Value = 4900
$ws.Range("M46").Value = -4712
$ws.Range("H55").Value = 1167.2142
$ws.Range("I55").Value = 440.33334
$ws.Range("K55").Value = 440.33334
$ws.Range("M55").Value = -267.33334
$ws.Range("H61").Value = 1337.2
$ws.Range("I61").Value = 1449.5
$ws.Range("K61").Value = 1449.5
$ws.Range("M61").Value = -1247.5
$ws.Range("H113").Value = 1337.2
$ws.Range("I113").Value = 1449.5
$ws.Range("K113").Value = 1449.5
$ws.Range("M113").Value = 720.5
$ws.Range("H136").Value = 3372.5
$ws.Range("I136").Value = 2498.3333
$ws.Range("J136").Value = 5995
$ws.Range("K136").Value = 7494.999899999999
$ws.Range("L136").Value = 17985
$ws.Range("M136").Value = -4944.999899999999
$ws.Range("N136").Value = -23085

# ---- Sheet WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H54").Value = 14878.477
$ws.Range("I54").Value = 10000
$ws.Range("J54").Value = 30489.6
$ws.Range("K54").Value = 10000
$ws.Range("L54").Value = 30489.6
$ws.Range("M54").Value = -9480
$ws.Range("N54").Value = -31529.6
